$wb = $excel.ActiveWorkbook

# --- Update version strings on the "compounds" sheet (sheet1) ---
# Row 3 (CIViC): source_version 23.08d -> 23.09d
# Row 2 (NCI Thesaurus): source_version 2023.06 -> 2023.09
# Update E3 first, then E2, so the newly created shared-string entries land in
# the same order as in the target workbook.
$wsCompounds = $wb.Worksheets.Item("compounds")

$civicVersion = $wsCompounds.Range("E3")
$civicVersion.NumberFormat = "@"
$civicVersion.Value = "23.09d"
$civicVersion.Style = "Normal"

$nciVersion = $wsCompounds.Range("E2")
$nciVersion.NumberFormat = "@"
$nciVersion.Value = "2023.09"
$nciVersion.Style = "Normal"

# --- Switch the active sheet from "biomarkers" back to "compounds" ---
$wsCompounds.Activate()
$wsCompounds.Range("E5").Select() | Out-Null
